$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "QuestionBlock"
$ws.Range("A2").Value = "Core Java"
$ws.Range("A3").Value = "Core Java"
$ws.Range("A4").Value = "Core Java"
$ws.Range("A5").Value = "Core Java"
$ws.Range("A6").Value = "Core Java"
$ws.Range("A7").Value = "Core Java"
$ws.Range("A8").Value = "Core Java"
$ws.Range("A9").Value = "Core Java"
$ws.Range("A10").Value = "Core Java"

$ws.Columns.Item(1).ColumnWidth = 11.42578125

$ws.Range("G16").Select()
